$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 4.4
$ws.Range("G2").Value = 5.8
$ws.Range("H2").Value = 1.77
$ws.Range("I2").Value = 1.96
$ws.Range("K2").Value = 4.7
$ws.Range("L2").Value = 1.31
$ws.Range("T2").Value = 1.75
$ws.Range("U2").Value = 2.12
$ws.Range("AC2").Value = 970
$ws.Range("AO2").Value = 1000
$ws.Range("F3").Value = 2.62
$ws.Range("K3").Value = 2.98
$ws.Range("L3").Value = 1.7
$ws.Range("V3").Value = 1.37
$ws.Range("F4").Value = 1.65
$ws.Range("G4").Value = 1.82
$ws.Range("I4").Value = 6.6
$ws.Range("J4").Value = 3.75
$ws.Range("L4").Value = 1.37
$ws.Range("N4").Value = 3.75
$ws.Range("O4").Value = 1.28
$ws.Range("P4").Value = 1.96
$ws.Range("Q4").Value = 1.84
$ws.Range("S4").Value = 3.15
$ws.Range("U4").Value = 1.98
$ws.Range("W4").Value = 2.2
$ws.Range("Y4").Value = 23
$ws.Range("Z4").Value = 60
$ws.Range("AD4").Value = 25
$ws.Range("AH4").Value = 22
$ws.Range("AL4").Value = 1000
$ws.Range("AN4").Value = 11
$ws.Range("F6").Value = 1.23
$ws.Range("H6").Value = 13
$ws.Range("I6").Value = 16.5
$ws.Range("J6").Value = 6.8
$ws.Range("K6").Value = 8.199999999999999
$ws.Range("L6").Value = 1.21
$ws.Range("N6").Value = 6.2
$ws.Range("O6").Value = 1.15
$ws.Range("P6").Value = 2.84
$ws.Range("Q6").Value = 1.44
$ws.Range("R6").Value = 1.73
$ws.Range("S6").Value = 2.1
$ws.Range("T6").Value = 1.97
$ws.Range("U6").Value = 1.84
$ws.Range("V6").Value = 1.06
$ws.Range("W6").Value = 4.5
$ws.Range("X6").Value = 40
$ws.Range("Y6").Value = 1000
$ws.Range("AC6").Value = 18
$ws.Range("AE6").Value = 250
$ws.Range("AG6").Value = 12
$ws.Range("AN6").Value = 3.8
$ws.Range("AO7").Value = 3.45
$ws.Range("I8").Value = 3.4
$ws.Range("U8").Value = 1.64
$ws.Range("V8").Value = 1.41
$ws.Range("Z8").Value = 21
$ws.Range("AC8").Value = 7
$ws.Range("AH8").Value = 29
$ws.Range("F9").Value = 1.74
$ws.Range("I9").Value = 7.4
$ws.Range("J9").Value = 3.4
$ws.Range("V9").Value = 1.18
$ws.Range("W9").Value = 2.22
$ws.Range("J10").Value = 3.4
$ws.Range("K10").Value = 3.9
$ws.Range("L10").Value = 1.48
$ws.Range("N10").Value = 3
$ws.Range("O10").Value = 1.41
$ws.Range("T10").Value = 2.04
$ws.Range("AB10").Value = 970
$ws.Range("AC10").Value = 1000
$ws.Range("O11").Value = 1.56
$ws.Range("X11").Value = 8.6
$ws.Range("Y11").Value = 11
$ws.Range("AB11").Value = 8.800000000000001
$ws.Range("AE11").Value = 60
$ws.Range("AF11").Value = 19
$ws.Range("AG11").Value = 14
$ws.Range("AH11").Value = 27
$ws.Range("AK11").Value = 42
$ws.Range("AN11").Value = 50
$ws.Range("G12").Value = 2.36
$ws.Range("K12").Value = 3.55
$ws.Range("O12").Value = 1.37
$ws.Range("S12").Value = 3.85
$ws.Range("G13").Value = 2.64
$ws.Range("I13").Value = 4.5
$ws.Range("J13").Value = 2.74
$ws.Range("W13").Value = 1.61
$ws.Range("AO13").Value = 95
$ws.Range("G14").Value = 1.7
$ws.Range("I14").Value = 7
$ws.Range("P14").Value = 1.88
$ws.Range("Q14").Value = 2
$ws.Range("W14").Value = 2.42
$ws.Range("Y14").Value = 970
$ws.Range("AH14").Value = 24
$ws.Range("AL14").Value = 40
$ws.Range("AM14").Value = 150
$ws.Range("AN14").Value = 11.5
$ws.Range("AB15").Value = 7.8
$ws.Range("AD15").Value = 44
$ws.Range("AD16").Value = 46
$ws.Range("AH16").Value = 42
